$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.303.08'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '1.873.86'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7091'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.07'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07782'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.83%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3105'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.06'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08417'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.30%  '
$ws.Range('D12').Value = '1.868.13'
$ws.Range('E12').Value = '  -0.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.240'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7160'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.11'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('D16').Value = '29.312.95'
$ws.Range('E16').Value = '  +0.08%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.086'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.44%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008282'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.52'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E20').Value = '  +0.51%  '
$ws.Range('D21').Value = '2.123.33'
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.752'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.64%  '
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1595'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.51'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.81%  '
$ws.Range('E27').Value = '  +0.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.50'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.09%  '
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.406'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.290'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.321'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05380'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.81%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.949'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.64%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.178'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.52%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7503'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.686'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('E38').Value = '  +1.26%  '
$ws.Range('D39').Value = '1.236.35'
$ws.Range('E39').Value = '  +6.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.730'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.69%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.473'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8922'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '72.39'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '108.81'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.00%  '
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('D46').Value = '2.017.11'
$ws.Range('E46').Value = '  -0.14%  '
$ws.Range('E47').Value = '  +9.97%  '
$ws.Range('E48').Value = '  -0.09%  '
$ws.Range('E49').Value = '  +0.14%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.445'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4337'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.79%  '
